$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for numeric-looking strings in Confidence % / Odds columns
$ws.Range("F2:G22").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 'Real Madrid v Juventus'
$ws.Range("B2").Value = 'Real Madrid'
$ws.Range("C2").Value = 'Champions League'
$ws.Range("D2").Value = '2025-10-22T19:00:00.000Z'
$ws.Range("E2").Value = '58/62 Win Tips'
$ws.Range("F2").Value = '94'
$ws.Range("G2").Value = '1.57'

# Row 3
$ws.Range("A3").Value = 'Bayern Munich v Club Brugge'
$ws.Range("B3").Value = 'Bayern Munich'
$ws.Range("C3").Value = 'Champions League'
$ws.Range("D3").Value = '2025-10-22T19:00:00.000Z'
$ws.Range("E3").Value = '56/56 Win Tips'
$ws.Range("F3").Value = '100'
$ws.Range("G3").Value = '1.22'

# Row 4
$ws.Range("A4").Value = 'Galatasaray v Bodo-Glimt'
$ws.Range("B4").Value = 'Galatasaray'
$ws.Range("C4").Value = 'Champions League'
$ws.Range("D4").Value = '2025-10-22T16:45:00.000Z'
$ws.Range("E4").Value = '49/61 Win Tips'
$ws.Range("F4").Value = '80'
$ws.Range("G4").Value = '1.60'

# Row 5
$ws.Range("A5").Value = 'Chelsea v Ajax'
$ws.Range("B5").Value = 'Chelsea'
$ws.Range("C5").Value = 'Champions League'
$ws.Range("D5").Value = '2025-10-22T19:00:00.000Z'
$ws.Range("E5").Value = '48/51 Win Tips'
$ws.Range("F5").Value = '94'
$ws.Range("G5").Value = '1.30'

# Row 6
$ws.Range("A6").Value = 'Atalanta v Slavia Prague'
$ws.Range("B6").Value = 'Atalanta'
$ws.Range("C6").Value = 'Champions League'
$ws.Range("D6").Value = '2025-10-22T19:00:00.000Z'
$ws.Range("E6").Value = '39/44 Win Tips'
$ws.Range("F6").Value = '89'
$ws.Range("G6").Value = '1.40'

# Row 7
$ws.Range("A7").Value = 'Athletic Bilbao v FK Qarabag'
$ws.Range("B7").Value = 'Athletic Bilbao'
$ws.Range("C7").Value = 'Champions League'
$ws.Range("D7").Value = '2025-10-22T16:45:00.000Z'
$ws.Range("E7").Value = '37/65 Win Tips'
$ws.Range("F7").Value = '57'
$ws.Range("G7").Value = '1.35'

# Row 8
$ws.Range("A8").Value = 'Sheff Wed v Middlesbrough'
$ws.Range("B8").Value = 'Middlesbrough'
$ws.Range("C8").Value = 'England Championship'
$ws.Range("D8").Value = '2025-10-22T19:00:00.000Z'
$ws.Range("E8").Value = '29/36 Win Tips'
$ws.Range("F8").Value = '81'
$ws.Range("G8").Value = '1.73'

# Row 9
$ws.Range("A9").Value = 'Wrexham v Oxford Utd'
$ws.Range("B9").Value = 'Wrexham'
$ws.Range("C9").Value = 'England Championship'
$ws.Range("D9").Value = '2025-10-22T18:45:00.000Z'
$ws.Range("E9").Value = '21/31 Win Tips'
$ws.Range("F9").Value = '68'
$ws.Range("G9").Value = '2.00'

# Row 10
$ws.Range("A10").Value = 'Eintracht Frankfurt v Liverpool'
$ws.Range("B10").Value = 'Liverpool'
$ws.Range("C10").Value = 'Champions League'
$ws.Range("D10").Value = '2025-10-22T19:00:00.000Z'
$ws.Range("E10").Value = '21/36 Win Tips'
$ws.Range("F10").Value = '58'
$ws.Range("G10").Value = '1.61'

# Row 11
$ws.Range("A11").Value = 'Swansea v QPR'
$ws.Range("B11").Value = 'Swansea'
$ws.Range("C11").Value = 'England Championship'
$ws.Range("D11").Value = '2025-10-22T18:45:00.000Z'
$ws.Range("E11").Value = '17/29 Win Tips'
$ws.Range("F11").Value = '59'
$ws.Range("G11").Value = '2.25'

# Row 12
$ws.Range("A12").Value = 'Monaco v Tottenham'
$ws.Range("B12").Value = 'Over 2.50'
$ws.Range("C12").Value = 'Champions League'
$ws.Range("D12").Value = '2025-10-22T19:00:00.000Z'
$ws.Range("E12").Value = '15/15 Win Tips'
$ws.Range("F12").Value = '100'
$ws.Range("G12").Value = '1.67'

# Row 13
$ws.Range("A13").Value = 'HJK Helsinki v SJK'
$ws.Range("B13").Value = 'Over 2.50'
$ws.Range("C13").Value = 'Finland Veikkausliiga'
$ws.Range("D13").Value = '2025-10-22T16:00:00.000Z'
$ws.Range("E13").Value = '14/14 Win Tips'
$ws.Range("F13").Value = '100'
$ws.Range("G13").Value = '1.40'

# Row 14
$ws.Range("A14").Value = 'Sporting v Marseille'
$ws.Range("B14").Value = 'Over 2.50'
$ws.Range("C14").Value = 'Champions League'
$ws.Range("D14").Value = '2025-10-22T19:00:00.000Z'
$ws.Range("E14").Value = '14/15 Win Tips'
$ws.Range("F14").Value = '93'
$ws.Range("G14").Value = '1.65'

# Row 15
$ws.Range("A15").Value = 'Chelsea v Sunderland'
$ws.Range("B15").Value = 'Chelsea'
$ws.Range("C15").Value = 'England Premier League'
$ws.Range("D15").Value = '2025-10-25T14:00:00.000Z'
$ws.Range("E15").Value = '14/18 Win Tips'
$ws.Range("F15").Value = '78'
$ws.Range("G15").Value = '1.42'

# Row 16
$ws.Range("A16").Value = 'Newcastle v Fulham'
$ws.Range("B16").Value = 'Newcastle'
$ws.Range("C16").Value = 'England Premier League'
$ws.Range("D16").Value = '2025-10-25T14:00:00.000Z'
$ws.Range("E16").Value = '13/15 Win Tips'
$ws.Range("F16").Value = '87'
$ws.Range("G16").Value = '1.62'

# Row 17
$ws.Range("A17").Value = 'Leeds v West Ham'
$ws.Range("B17").Value = 'Leeds'
$ws.Range("C17").Value = 'England Premier League'
$ws.Range("D17").Value = '2025-10-24T19:00:00.000Z'
$ws.Range("E17").Value = '11/15 Win Tips'
$ws.Range("F17").Value = '73'
$ws.Range("G17").Value = '2.05'

# Row 18
$ws.Range("A18").Value = 'Arsenal v Crystal Palace'
$ws.Range("B18").Value = 'Arsenal'
$ws.Range("C18").Value = 'England Premier League'
$ws.Range("D18").Value = '2025-10-26T14:00:00.000Z'
$ws.Range("E18").Value = '11/13 Win Tips'
$ws.Range("F18").Value = '85'
$ws.Range("G18").Value = '1.48'

# Row 19
$ws.Range("A19").Value = 'Watford v West Brom'
$ws.Range("B19").Value = 'Watford'
$ws.Range("C19").Value = 'England Championship'
$ws.Range("D19").Value = '2025-10-22T18:45:00.000Z'
$ws.Range("E19").Value = '10/21 Win Tips'
$ws.Range("F19").Value = '48'
$ws.Range("G19").Value = '2.90'

# Row 20
$ws.Range("A20").Value = 'Brentford v Liverpool'
$ws.Range("B20").Value = 'Liverpool'
$ws.Range("C20").Value = 'England Premier League'
$ws.Range("D20").Value = '2025-10-25T19:00:00.000Z'
$ws.Range("E20").Value = '10/15 Win Tips'
$ws.Range("F20").Value = '67'
$ws.Range("G20").Value = '1.75'

# Row 21
$ws.Range("A21").Value = 'Bohemians 1905 v Mlada Boleslav'
$ws.Range("B21").Value = 'Bohemians 1905'
$ws.Range("C21").Value = 'Czech Republic First League'
$ws.Range("D21").Value = '2025-10-22T16:30:00.000Z'
$ws.Range("E21").Value = '9/14 Win Tips'
$ws.Range("F21").Value = '64'
$ws.Range("G21").Value = '2.05'

# Row 22
$ws.Range("A22").Value = 'UEFA Champions League 2025-26'
$ws.Range("B22").Value = 'PSG'
$ws.Range("C22").Value = 'Champions League'
$ws.Range("D22").Value = '2026-05-30T19:00:00.000Z'
$ws.Range("E22").Value = '4/25 Win Tips'
$ws.Range("F22").Value = '16'
$ws.Range("G22").Value = '5.50'
